$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("miasta")

# Row 2 held "Gdańsk " with wage 10017.25 -> change city to "Radom" and wage to 7368
$ws.Range("A2").Value = "Radom"
$ws.Range("B2").Value = 7368

# Copy formatting from a neighboring already-styled city cell (A3) onto A2,
# then add a new styled (but empty) row 8 below the existing data, matching A3's style.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B3").Select() | Out-Null
